# Weekly refresh of "Fruta, Vega Central Mapocho de Santiago - Coco" data:
# columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen) and
# S (Precio $/Kg) are reshuffled/updated per row for rows 2-41.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44363; $ws.Range("L2").Value2 = "Primera"; $ws.Range("M2").Value2 = 150; $ws.Range("N2").Value2 = 21000; $ws.Range("O2").Value2 = 22000; $ws.Range("P2").Value2 = 21500; $ws.Range("R2").Value2 = "Perú"; $ws.Range("S2").Value2 = 1075
$ws.Range("D3").Value2 = 44166; $ws.Range("L3").Value2 = "Primera"; $ws.Range("M3").Value2 = 120; $ws.Range("N3").Value2 = 28000; $ws.Range("O3").Value2 = 28000; $ws.Range("P3").Value2 = 28000; $ws.Range("R3").Value2 = "Perú"; $ws.Range("S3").Value2 = 1400
$ws.Range("D4").Value2 = 44438; $ws.Range("L4").Value2 = "Primera"; $ws.Range("M4").Value2 = 25; $ws.Range("N4").Value2 = 21000; $ws.Range("O4").Value2 = 21000; $ws.Range("P4").Value2 = 21000; $ws.Range("R4").Value2 = "Perú"; $ws.Range("S4").Value2 = 1050
$ws.Range("D5").Value2 = 44356; $ws.Range("L5").Value2 = "Primera"; $ws.Range("M5").Value2 = 100; $ws.Range("N5").Value2 = 20000; $ws.Range("O5").Value2 = 21000; $ws.Range("P5").Value2 = 20500; $ws.Range("R5").Value2 = "Perú"; $ws.Range("S5").Value2 = 1025
$ws.Range("D6").Value2 = 44165; $ws.Range("L6").Value2 = "Primera"; $ws.Range("M6").Value2 = 300; $ws.Range("N6").Value2 = 27000; $ws.Range("O6").Value2 = 28000; $ws.Range("P6").Value2 = 27500; $ws.Range("R6").Value2 = "Perú"; $ws.Range("S6").Value2 = 1375
$ws.Range("D7").Value2 = 44277; $ws.Range("L7").Value2 = "Primera"; $ws.Range("M7").Value2 = 60; $ws.Range("N7").Value2 = 24000; $ws.Range("O7").Value2 = 24000; $ws.Range("P7").Value2 = 24000; $ws.Range("R7").Value2 = "Perú"; $ws.Range("S7").Value2 = 1200
$ws.Range("D8").Value2 = 44299; $ws.Range("L8").Value2 = "Primera"; $ws.Range("M8").Value2 = 150; $ws.Range("N8").Value2 = 19000; $ws.Range("O8").Value2 = 20000; $ws.Range("P8").Value2 = 19500; $ws.Range("R8").Value2 = "Perú"; $ws.Range("S8").Value2 = 975
$ws.Range("D9").Value2 = 44326; $ws.Range("L9").Value2 = "Primera"; $ws.Range("M9").Value2 = 40; $ws.Range("N9").Value2 = 22000; $ws.Range("O9").Value2 = 22000; $ws.Range("P9").Value2 = 22000; $ws.Range("R9").Value2 = "Perú"; $ws.Range("S9").Value2 = 1100
$ws.Range("D10").Value2 = 44300; $ws.Range("L10").Value2 = "Primera"; $ws.Range("M10").Value2 = 150; $ws.Range("N10").Value2 = 19000; $ws.Range("O10").Value2 = 20000; $ws.Range("P10").Value2 = 19500; $ws.Range("R10").Value2 = "Perú"; $ws.Range("S10").Value2 = 975
$ws.Range("D11").Value2 = 44270; $ws.Range("L11").Value2 = "Primera"; $ws.Range("M11").Value2 = 50; $ws.Range("N11").Value2 = 24000; $ws.Range("O11").Value2 = 24000; $ws.Range("P11").Value2 = 24000; $ws.Range("R11").Value2 = "Perú"; $ws.Range("S11").Value2 = 1200
$ws.Range("D12").Value2 = 44357; $ws.Range("L12").Value2 = "Primera"; $ws.Range("M12").Value2 = 200; $ws.Range("N12").Value2 = 20000; $ws.Range("O12").Value2 = 21000; $ws.Range("P12").Value2 = 20500; $ws.Range("R12").Value2 = "Perú"; $ws.Range("S12").Value2 = 1025
$ws.Range("D13").Value2 = 44312; $ws.Range("L13").Value2 = "Primera"; $ws.Range("M13").Value2 = 50; $ws.Range("N13").Value2 = 22000; $ws.Range("O13").Value2 = 22000; $ws.Range("P13").Value2 = 22000; $ws.Range("R13").Value2 = "Perú"; $ws.Range("S13").Value2 = 1100
$ws.Range("D14").Value2 = 44529; $ws.Range("L14").Value2 = "Primera"; $ws.Range("M14").Value2 = 34; $ws.Range("N14").Value2 = 28000; $ws.Range("O14").Value2 = 28000; $ws.Range("P14").Value2 = 28000; $ws.Range("R14").Value2 = "Perú"; $ws.Range("S14").Value2 = 1400
$ws.Range("D15").Value2 = 44396; $ws.Range("L15").Value2 = "Primera"; $ws.Range("M15").Value2 = 45; $ws.Range("N15").Value2 = 22000; $ws.Range("O15").Value2 = 22000; $ws.Range("P15").Value2 = 22000; $ws.Range("R15").Value2 = "Perú"; $ws.Range("S15").Value2 = 1100
$ws.Range("D16").Value2 = 44305; $ws.Range("L16").Value2 = "Primera"; $ws.Range("M16").Value2 = 40; $ws.Range("N16").Value2 = 24000; $ws.Range("O16").Value2 = 24000; $ws.Range("P16").Value2 = 24000; $ws.Range("R16").Value2 = "Perú"; $ws.Range("S16").Value2 = 1200
$ws.Range("D17").Value2 = 44302; $ws.Range("L17").Value2 = "Primera"; $ws.Range("M17").Value2 = 100; $ws.Range("N17").Value2 = 19000; $ws.Range("O17").Value2 = 20000; $ws.Range("P17").Value2 = 19500; $ws.Range("R17").Value2 = "Perú"; $ws.Range("S17").Value2 = 975
$ws.Range("D18").Value2 = 44522; $ws.Range("L18").Value2 = "Primera"; $ws.Range("M18").Value2 = 25; $ws.Range("N18").Value2 = 30000; $ws.Range("O18").Value2 = 30000; $ws.Range("P18").Value2 = 30000; $ws.Range("R18").Value2 = "Perú"; $ws.Range("S18").Value2 = 1500
$ws.Range("D19").Value2 = 44424; $ws.Range("L19").Value2 = "Primera"; $ws.Range("M19").Value2 = 70; $ws.Range("N19").Value2 = 24000; $ws.Range("O19").Value2 = 25000; $ws.Range("P19").Value2 = 24429; $ws.Range("R19").Value2 = "Perú"; $ws.Range("S19").Value2 = 1221
$ws.Range("D20").Value2 = 44445; $ws.Range("L20").Value2 = "Primera"; $ws.Range("M20").Value2 = 35; $ws.Range("N20").Value2 = 20000; $ws.Range("O20").Value2 = 20000; $ws.Range("P20").Value2 = 20000; $ws.Range("R20").Value2 = "Perú"; $ws.Range("S20").Value2 = 1000
$ws.Range("D21").Value2 = 44452; $ws.Range("L21").Value2 = "Primera"; $ws.Range("M21").Value2 = 35; $ws.Range("N21").Value2 = 21000; $ws.Range("O21").Value2 = 22000; $ws.Range("P21").Value2 = 21429; $ws.Range("R21").Value2 = "Perú"; $ws.Range("S21").Value2 = 1071
$ws.Range("D22").Value2 = 44435; $ws.Range("L22").Value2 = "Primera"; $ws.Range("M22").Value2 = 60; $ws.Range("N22").Value2 = 25000; $ws.Range("O22").Value2 = 25000; $ws.Range("P22").Value2 = 25000; $ws.Range("R22").Value2 = "Perú"; $ws.Range("S22").Value2 = 1250
$ws.Range("D23").Value2 = 44613; $ws.Range("L23").Value2 = "Primera"; $ws.Range("M23").Value2 = 60; $ws.Range("N23").Value2 = 30000; $ws.Range("O23").Value2 = 30000; $ws.Range("P23").Value2 = 30000; $ws.Range("R23").Value2 = "Perú"; $ws.Range("S23").Value2 = 1500
$ws.Range("D24").Value2 = 44372; $ws.Range("L24").Value2 = "Primera"; $ws.Range("M24").Value2 = 60; $ws.Range("N24").Value2 = 20000; $ws.Range("O24").Value2 = 21000; $ws.Range("P24").Value2 = 20667; $ws.Range("R24").Value2 = "Perú"; $ws.Range("S24").Value2 = 1033
$ws.Range("D25").Value2 = 44620; $ws.Range("L25").Value2 = "Primera"; $ws.Range("M25").Value2 = 60; $ws.Range("N25").Value2 = 22000; $ws.Range("O25").Value2 = 22000; $ws.Range("P25").Value2 = 22000; $ws.Range("R25").Value2 = "Perú"; $ws.Range("S25").Value2 = 1100
$ws.Range("D26").Value2 = 44830; $ws.Range("L26").Value2 = "Primera"; $ws.Range("M26").Value2 = 200; $ws.Range("N26").Value2 = 30000; $ws.Range("O26").Value2 = 30000; $ws.Range("P26").Value2 = 30000; $ws.Range("R26").Value2 = "Perú"; $ws.Range("S26").Value2 = 1500
$ws.Range("D27").Value2 = 44442; $ws.Range("L27").Value2 = "Primera"; $ws.Range("M27").Value2 = 30; $ws.Range("N27").Value2 = 22000; $ws.Range("O27").Value2 = 22000; $ws.Range("P27").Value2 = 22000; $ws.Range("R27").Value2 = "Perú"; $ws.Range("S27").Value2 = 1100
$ws.Range("D28").Value2 = 44760; $ws.Range("L28").Value2 = "Primera"; $ws.Range("M28").Value2 = 300; $ws.Range("N28").Value2 = 24000; $ws.Range("O28").Value2 = 25000; $ws.Range("P28").Value2 = 24500; $ws.Range("R28").Value2 = "Perú"; $ws.Range("S28").Value2 = 1225
$ws.Range("D29").Value2 = 44354; $ws.Range("L29").Value2 = "Primera"; $ws.Range("M29").Value2 = 150; $ws.Range("N29").Value2 = 21000; $ws.Range("O29").Value2 = 22000; $ws.Range("P29").Value2 = 21500; $ws.Range("R29").Value2 = "Perú"; $ws.Range("S29").Value2 = 1075
$ws.Range("D30").Value2 = 44355; $ws.Range("L30").Value2 = "Primera"; $ws.Range("M30").Value2 = 200; $ws.Range("N30").Value2 = 20000; $ws.Range("O30").Value2 = 21000; $ws.Range("P30").Value2 = 20500; $ws.Range("R30").Value2 = "Ecuador"; $ws.Range("S30").Value2 = 1025
$ws.Range("D31").Value2 = 44473; $ws.Range("L31").Value2 = "Primera"; $ws.Range("M31").Value2 = 40; $ws.Range("N31").Value2 = 24000; $ws.Range("O31").Value2 = 24000; $ws.Range("P31").Value2 = 24000; $ws.Range("R31").Value2 = "Perú"; $ws.Range("S31").Value2 = 1200
$ws.Range("D32").Value2 = 44382; $ws.Range("L32").Value2 = "Primera"; $ws.Range("M32").Value2 = 200; $ws.Range("N32").Value2 = 19000; $ws.Range("O32").Value2 = 20000; $ws.Range("P32").Value2 = 19500; $ws.Range("R32").Value2 = "Perú"; $ws.Range("S32").Value2 = 975
$ws.Range("D33").Value2 = 44263; $ws.Range("L33").Value2 = "Segunda"; $ws.Range("M33").Value2 = 150; $ws.Range("N33").Value2 = 15000; $ws.Range("O33").Value2 = 15000; $ws.Range("P33").Value2 = 15000; $ws.Range("R33").Value2 = "Perú"; $ws.Range("S33").Value2 = 750
$ws.Range("D34").Value2 = 44350; $ws.Range("L34").Value2 = "Primera"; $ws.Range("M34").Value2 = 90; $ws.Range("N34").Value2 = 21000; $ws.Range("O34").Value2 = 22000; $ws.Range("P34").Value2 = 21556; $ws.Range("R34").Value2 = "Perú"; $ws.Range("S34").Value2 = 1078
$ws.Range("D35").Value2 = 44298; $ws.Range("L35").Value2 = "Primera"; $ws.Range("M35").Value2 = 240; $ws.Range("N35").Value2 = 19000; $ws.Range("O35").Value2 = 20000; $ws.Range("P35").Value2 = 19500; $ws.Range("R35").Value2 = "Perú"; $ws.Range("S35").Value2 = 975
$ws.Range("D36").Value2 = 44284; $ws.Range("L36").Value2 = "Primera"; $ws.Range("M36").Value2 = 40; $ws.Range("N36").Value2 = 23000; $ws.Range("O36").Value2 = 23000; $ws.Range("P36").Value2 = 23000; $ws.Range("R36").Value2 = "Perú"; $ws.Range("S36").Value2 = 1150
$ws.Range("D37").Value2 = 44333; $ws.Range("L37").Value2 = "Primera"; $ws.Range("M37").Value2 = 30; $ws.Range("N37").Value2 = 22000; $ws.Range("O37").Value2 = 22000; $ws.Range("P37").Value2 = 22000; $ws.Range("R37").Value2 = "Perú"; $ws.Range("S37").Value2 = 1100
$ws.Range("D38").Value2 = 44417; $ws.Range("L38").Value2 = "Primera"; $ws.Range("M38").Value2 = 30; $ws.Range("N38").Value2 = 24000; $ws.Range("O38").Value2 = 24000; $ws.Range("P38").Value2 = 24000; $ws.Range("R38").Value2 = "Perú"; $ws.Range("S38").Value2 = 1200
$ws.Range("D39").Value2 = 44410; $ws.Range("L39").Value2 = "Primera"; $ws.Range("M39").Value2 = 40; $ws.Range("N39").Value2 = 25000; $ws.Range("O39").Value2 = 25000; $ws.Range("P39").Value2 = 25000; $ws.Range("R39").Value2 = "Perú"; $ws.Range("S39").Value2 = 1250
$ws.Range("D40").Value2 = 44431; $ws.Range("L40").Value2 = "Primera"; $ws.Range("M40").Value2 = 60; $ws.Range("N40").Value2 = 25000; $ws.Range("O40").Value2 = 25000; $ws.Range("P40").Value2 = 25000; $ws.Range("R40").Value2 = "Perú"; $ws.Range("S40").Value2 = 1250
$ws.Range("D41").Value2 = 44365; $ws.Range("L41").Value2 = "Primera"; $ws.Range("M41").Value2 = 150; $ws.Range("N41").Value2 = 20000; $ws.Range("O41").Value2 = 21000; $ws.Range("P41").Value2 = 20500; $ws.Range("R41").Value2 = "Perú"; $ws.Range("S41").Value2 = 1025
